# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" positioned right before the
#    "总计" (totals) sheet, and populate it with the per-fund holdings
#    data for the 2022-Q1 quarter (same layout as the other quarterly
#    sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# 2. Insert a new summary row at the top of the "总计" sheet's data for
#    2022-Q1 (持有数量=8, 持有市值=5.21), pushing the existing rows down
#    and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$totalsBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalsBefore)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous

# Per-fund holdings data. Columns D/E/F/G are kept as *text* (matching the
# other quarterly sheets) even though they look numeric, by prefixing the
# literal with an apostrophe; A (index) and H (rank) are real numbers.
$rows = @(
    @("510810", "汇添富中证上海国企ETF", "68.43", "99.71", "7.33", "5.0159", 4),
    @("009073", "德邦惠利混合A", "2.56", "32.80", "1.76", "0.0451", 5),
    @("519616", "银河君信灵活配置混合A", "4.54", "24.42", "0.84", "0.0381", 6),
    @("519618", "银河君信灵活配置混合I", "4.54", "24.42", "0.84", "0.0381", 6),
    @("000433", "安信鑫发优选混合", "1.23", "67.20", "2.55", "0.0314", 5),
    @("009074", "德邦惠利混合C", "1.01", "32.80", "1.76", "0.0178", 5),
    @("004260", "德邦稳盈增长灵活配置混合", "0.83", "27.99", "1.71", "0.0142", 5),
    @("519617", "银河君信灵活配置混合C", "0.64", "24.42", "0.84", "0.0054", 6)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q1.Range("A$r")
    $aCell.Value = $r - 2
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $q1.Range("B$r").Value = "'" + $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = "'" + $row[2]
    $q1.Range("E$r").Value = "'" + $row[3]
    $q1.Range("F$r").Value = "'" + $row[4]
    $q1.Range("G$r").Value = "'" + $row[5]
    $q1.Range("H$r").Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------------
# NOTE: worksheet proxies in this runtime re-resolve by *position*, so the
# handle obtained before the insert above now points at the new "2022-Q1"
# sheet (which took over slot 6). Re-fetch "总计" by name now that it has
# settled into its final slot (7).
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# The inserted row inherits the header row's bold/centered formatting;
# reset the plain data cells (B:D) back to the default look used by every
# other data row, then (re)apply the special index-column style to A2.
$totals.Range("B2:D2").ClearFormats()

$a2 = $totals.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 8
$totals.Range("D2").Value = 5.21

# Renumber the index column for the rows that got pushed down (they kept
# their old 0-based values after the insert).
for ($row = 3; $row -le 7; $row++) {
    $totals.Range("A$row").Value = $row - 2
}

# Restore the originally active sheet/selection (Worksheets.Add above
# activates the freshly inserted sheet as a side effect, like real Excel).
$wb.Worksheets.Item(1).Activate() | Out-Null
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null

